# Super Meat Boy presentation - 26-5-21 update
#  1. Reorder the "Personajes" slide so it comes before "Escenarios" and
#     "Mecanicas principales" (was: Escenarios, Mecanicas principales,
#     Personajes -> now: Personajes, Escenarios, Mecanicas principales).
#  2. Match that same new order on the agenda ("INDICE") slide.
#  3. Small typo fix on the "Trabajo pendiente" slide.

$p = $ppt.ActivePresentation

# --- 1. Reorder slides: move the "Personajes" slide (currently #6) so it
#        lands right before the "Escenarios" slide (currently #4). ---
$personajesSlide = $p.Slides.Item(6)
$personajesSlide.MoveTo(4)

# --- 2. Fix the agenda slide (slide 2, "INDICE") bullet order so
#        "Personajes" is listed before "Escenarios" / "Mecanicas
#        principales", matching the new slide order. ---
$indice = $p.Slides.Item(2)
$body = $indice.Shapes.Item(2).TextFrame.TextRange

# Paragraph 2 = "Escenarios", 3 = "Mecanicas principales", 4 = "Personajes"
# (in that order) before the edit; retarget them to "Personajes",
# "Escenarios", "Mecanicas principales". Writing a throwaway value first
# avoids the COM text-diff engine treating the new/old strings as sharing
# a common prefix (which would otherwise split the run in two).
$para = $body.Paragraphs(2, 1)
$para.Text = " "
$para = $body.Paragraphs(2, 1)
$para.Text = "Personajes"

$para = $body.Paragraphs(3, 1)
$para.Text = " "
$para = $body.Paragraphs(3, 1)
$para.Text = "Escenarios"

$para = $body.Paragraphs(4, 1)
$para.Text = " "
$para = $body.Paragraphs(4, 1)
$para.Text = "Mecánicas principales"

# --- 3. Typo fix on the "Trabajo pendiente" slide (slide 8): remove the
#        stray "de" in "Diferentes de animaciones ." ---
$pendiente = $p.Slides.Item(8)
$bullets = $pendiente.Shapes.Item(2).TextFrame.TextRange
$hit = $bullets.Find("Diferentes de animaciones .", 0, $false, $false)
if ($hit) {
    $hit.Text = "Diferentes animaciones ."
}
